$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data set for rows 2-10 (9 rows total), columns A-T.
# Column order: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E..T = numeric metrics.
$data = @(
    @("ECs", "Myoc", "Fzd4", "ECs", 2, 0.6666666666666666, 0.1664173333333333, 0.499252, 0.009928329481286188, 0.01026989293949489, 3, 1, 29.22839533333334, 87.685186, 0.4452121315669252, 0.4988525987881078, 4.864111608985778, 43.777004480872, 0.004420212731262169, 0.005123162782142665),
    @("ECs", "Myoc", "Fzd4", "FAPs", 2, 0.6666666666666666, 0.1664173333333333, 0.499252, 0.009928329481286188, 0.01026989293949489, 3, 1, 15.24435933333334, 45.73307800000001, 0.232204800700274, 0.2601815181287206, 2.536925628628445, 22.832330657656, 0.002305405768488714, 0.00267203633601721),
    @("ECs", "Myoc", "Fzd4", "sCs", 2, 0.6666666666666666, 0.1664173333333333, 0.499252, 0.009928329481286188, 0.01026989293949489, 2, 1, 21.177737, 42.355474, 0.3225830677328006, 0.2409658830831714, 3.524342517574667, 21.146055105448, 0.003202710981535303, 0.002474693821335012),
    @("FAPs", "Myoc", "Fzd4", "ECs", 3, 1, 14.92301466666667, 44.76904399999999, 0.8902955208876449, 0.920924280490686, 3, 1, 29.22839533333334, 87.685186, 0.4452121315669252, 0.4988525987881078, 436.1757722424649, 3925.581950182183, 0.3963703665788744, 0.459405470609847),
    @("FAPs", "Myoc", "Fzd4", "FAPs", 3, 1, 14.92301466666667, 44.76904399999999, 0.8902955208876449, 0.920924280490686, 3, 1, 15.24435933333334, 45.73307800000001, 0.232204800700274, 0.2601815181287206, 227.4917979152702, 2047.426181237432, 0.2067308939920622, 0.2396074773796664),
    @("FAPs", "Myoc", "Fzd4", "sCs", 3, 1, 14.92301466666667, 44.76904399999999, 0.8902955208876449, 0.920924280490686, 2, 1, 21.177737, 42.355474, 0.3225830677328006, 0.2409658830831714, 316.0356798578093, 1896.214079146856, 0.2871942603167081, 0.2219113325011723),
    @("sCs", "Myoc", "Fzd4", "ECs", 2, 1, 1.6724345, 3.344869, 0.09977614963106883, 0.06880582656981911, 3, 1, 29.22839533333334, 87.685186, 0.4452121315669252, 0.4988525987881078, 48.88257673510567, 293.295460410634, 0.04442155225678864, 0.0343239653961181),
    @("sCs", "Myoc", "Fzd4", "FAPs", 2, 1, 1.6724345, 3.344869, 0.09977614963106883, 0.06880582656981911, 3, 1, 15.24435933333334, 45.73307800000001, 0.232204800700274, 0.2601815181287206, 25.49519247946367, 152.971154876782, 0.02316850093972305, 0.017902004413037),
    @("sCs", "Myoc", "Fzd4", "sCs", 2, 1, 1.6724345, 3.344869, 0.09977614963106883, 0.06880582656981911, 2, 1, 21.177737, 42.355474, 0.3225830677328006, 0.2409658830831714, 35.4183779907265, 141.673511962906, 0.03218609643455712, 0.016579856760664)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $rowValues[$c]
    }
}

# Clear any leftover cells beyond the new data range (rows 11+ in case of stray content)
$ws.Range("A11:T20").ClearContents()
